$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "27.321.01"
Set-TextValue 2 5 "  -0.77%  "
Set-TextValue 3 4 "1.710.72"
Set-TextValue 3 5 "  -0.66%  "
Set-TextValue 4 5 "  +0.01%  "
Set-TextValue 5 4 "224.28"
Set-TextValue 5 5 "  -0.65%  "
Set-TextValue 6 4 "0.5298"
Set-TextValue 6 5 "  -1.33%  "
Set-TextValue 7 5 "  +0.01%  "
Set-TextValue 8 4 "0.06702"
Set-TextValue 8 5 "  +1.53%  "
Set-TextValue 9 4 "0.2662"
Set-TextValue 10 4 "20.84"
Set-TextValue 10 5 "  -3.84%  "
Set-TextValue 11 4 "0.07676"
Set-TextValue 11 5 "  -0.59%  "
Set-TextValue 12 4 "4.508"
Set-TextValue 12 5 "  -2.36%  "
Set-TextValue 13 4 "1.946.55"
Set-TextValue 13 5 "  -0.66%  "
Set-TextValue 14 4 "1.710.06"
Set-TextValue 14 5 "  -0.64%  "
Set-TextValue 15 4 "0.5818"
Set-TextValue 15 5 "  -0.30%  "
Set-TextValue 16 4 "0.0₅8236"
Set-TextValue 16 5 "  -0.83%  "
Set-TextValue 17 4 "68.19"
Set-TextValue 17 5 "  +0.49%  "
Set-TextValue 18 4 "27.333.91"
Set-TextValue 19 4 "226.31"
Set-TextValue 19 5 "  +2.83%  "
Set-TextValue 20 5 "  -0.06%  "
Set-TextValue 21 4 "4.627"
Set-TextValue 21 5 "  -2.04%  "
Set-TextValue 22 5 "  -2.36%  "
Set-TextValue 23 4 "6.005"
Set-TextValue 23 5 "  -1.14%  "
Set-TextValue 24 4 "1.004"
Set-TextValue 24 5 "  -0.03%  "
Set-TextValue 25 5 "  -2.54%  "
Set-TextValue 26 4 "1.696"
Set-TextValue 26 5 "  -2.28%  "
Set-TextValue 27 4 "0.1206"
Set-TextValue 27 5 "  -2.47%  "
Set-TextValue 28 4 "7.233"
Set-TextValue 28 5 "  -2.25%  "
Set-TextValue 29 4 "16.30"
Set-TextValue 29 5 "  -1.77%  "
Set-TextValue 30 4 "0.05364"
Set-TextValue 30 5 "  -3.78%  "
Set-TextValue 31 4 "1.290"
Set-TextValue 31 5 "  -0.83%  "
Set-TextValue 32 5 "  -2.11%  "
Set-TextValue 33 4 "3.433"
Set-TextValue 34 5 "  -1.62%  "
Set-TextValue 35 4 "2.874"
Set-TextValue 35 5 "  +1.21%  "
Set-TextValue 36 4 "0.9500"
Set-TextValue 36 5 "  -1.27%  "
Set-TextValue 37 4 "2.394"
Set-TextValue 37 5 "  -1.44%  "
Set-TextValue 38 4 "0.5836"
Set-TextValue 38 5 "  -2.17%  "
Set-TextValue 39 4 "0.01634"
Set-TextValue 39 5 "  -0.86%  "
Set-TextValue 40 4 "1.080.43"
Set-TextValue 40 5 "  +2.58%  "
Set-TextValue 41 4 "5.785"
Set-TextValue 41 5 "  -2.38%  "
Set-TextValue 42 4 "1.004"
Set-TextValue 42 5 "  +0.02%  "
Set-TextValue 43 4 "0.8392"
Set-TextValue 43 5 "  -1.90%  "
Set-TextValue 44 4 "100.74"
Set-TextValue 44 5 "  -0.50%  "
Set-TextValue 45 4 "1.853.83"
Set-TextValue 45 5 "  -0.68%  "
Set-TextValue 46 4 "0.0₈116"
Set-TextValue 46 5 "  +0.87%  "
Set-TextValue 47 5 "  -2.38%  "
Set-TextValue 48 4 "0.4529"
Set-TextValue 48 5 "  +2.14%  "
Set-TextValue 49 4 "1.003"
Set-TextValue 49 5 "  +0.25%  "
Set-TextValue 50 4 "8.076"
Set-TextValue 50 5 "  -1.75%  "
Set-TextValue 51 4 "0.05232"
Set-TextValue 51 5 "  -0.23%  "
